$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.3351768771135573
$ws.Range("B2").Value = 0.5060638175655189
$ws.Range("C2").Value = 0.1796546986643375
$ws.Range("D2").Value = 0.4195427544048741
$ws.Range("E2").Value = 0.245578189146106
$ws.Range("K2").Value = 2.346238139794901
$ws.Range("L2").Value = 3.542446722958632
$ws.Range("M2").Value = 1.257582890650363
$ws.Range("N2").Value = 2.936799280834118
$ws.Range("O2").Value = 1.719047324022742
$ws.Range("P2").Value = 66.71236
$ws.Range("Q2").Value = 158.8199294437367
$ws.Range("R2").Value = -0.9222002297039863
$ws.Range("S2").Value = 93.94103923246006
$ws.Range("T2").Value = 41.68995937282172
$ws.Range("U2").Value = 0.2013905106781432
$ws.Range("V2").Value = 0.3608144527543763
$ws.Range("W2").Value = 0.01761398589523827
$ws.Range("X2").Value = 0.2937437689855742
$ws.Range("Y2").Value = 0.09906292406488931
$ws.Range("Z2").Value = 0.806474547716027
$ws.Range("AA2").Value = 1.013577409345779
$ws.Range("AB2").Value = 0.2791071634945997
$ws.Range("AC2").Value = 0.9052320647278638
$ws.Range("AD2").Value = 0.6232276245013814
$ws.Range("F3").Value = 7.050185388242457
$ws.Range("G3").Value = 10.39984920659461
$ws.Range("H3").Value = 3.86079638696003
$ws.Range("I3").Value = 8.707853521770621
$ws.Range("J3").Value = 5.299282209844216
$ws.Range("K3").Value = 2.326561178120011
$ws.Range("L3").Value = 3.431950238176221
$ws.Range("M3").Value = 1.27406280769681
$ws.Range("N3").Value = 2.873591662184305
$ws.Range("O3").Value = 1.748763129248591
$ws.Range("P3").Value = 66.78951000000001
$ws.Range("Q3").Value = 114.5183013445008
$ws.Range("R3").Value = 52.14962197739374
$ws.Range("S3").Value = 75.60569506759896
$ws.Range("T3").Value = 56.39191629479538
$ws.Range("U3").Value = 0.2001324207050712
$ws.Range("V3").Value = 0.3491814829917947
$ws.Range("W3").Value = 0.0295351232119927
$ws.Range("X3").Value = 0.2859781489362377
$ws.Range("Y3").Value = 0.1052450444196071
$ws.Range("Z3").Value = 0.8198098386603069
$ws.Range("AA3").Value = 0.9935588253029293
$ws.Range("AB3").Value = 0.4174250330544992
$ws.Range("AC3").Value = 0.9111008746823749
$ws.Range("AD3").Value = 0.6677358467506402
$ws.Range("A4").Value = 0.3350558419577123
$ws.Range("B4").Value = 0.5086499483008653
$ws.Range("C4").Value = 0.1800329909208327
$ws.Range("D4").Value = 0.4191459405920589
$ws.Range("E4").Value = 0.2454954139778766
$ws.Range("F4").Value = 7.053464883417598
$ws.Range("G4").Value = 10.4033846660129
$ws.Range("H4").Value = 3.861647242603307
$ws.Range("I4").Value = 8.702338044188473
$ws.Range("J4").Value = 5.317217194305628
$ws.Range("K4").Value = 2.363622771069657
$ws.Range("L4").Value = 4.233048760917629
$ws.Range("M4").Value = 1.025475250545921
$ws.Range("N4").Value = 3.165455211722798
$ws.Range("O4").Value = 1.518695666687967
$ws.Range("P4").Value = 66.27011
$ws.Range("Q4").Value = 163.82047538298
$ws.Range("R4").Value = -0.6129810202759727
$ws.Range("S4").Value = 96.23347761342018
$ws.Range("T4").Value = 39.78443561287411
$ws.Range("U4").Value = 0.1986317842501169
$ws.Range("V4").Value = 0.4167168209495163
$ws.Range("W4").Value = -0.009228800358541341
$ws.Range("X4").Value = 0.3265265691653199
$ws.Range("Y4").Value = 0.06783025487022593
$ws.Range("Z4").Value = 0.7663562657945971
$ws.Range("AA4").Value = 1.025878516857799
$ws.Range("AB4").Value = 0.08141683977251001
$ws.Range("AC4").Value = 0.9214637787715185
$ws.Range("AD4").Value = 0.500462807955672
